$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Cells.Item(112, 8).Value = 1933.4572
$ws.Cells.Item(112, 9).Value = 0
$ws.Cells.Item(112, 10).Value = 1933.4572
$ws.Cells.Item(112, 11).Value = 0
$ws.Cells.Item(112, 12).Value = 5800.3716
$ws.Cells.Item(112, 13).ClearContents()
$ws.Cells.Item(112, 14).Value = -8016.3716
# Row 129
$ws.Cells.Item(129, 8).Value = 1038.7894
$ws.Cells.Item(129, 9).Value = 649.25
$ws.Cells.Item(129, 10).Value = 1084.6177
$ws.Cells.Item(129, 11).Value = 1947.75
$ws.Cells.Item(129, 12).Value = 3253.8531
$ws.Cells.Item(129, 13).Value = 3052.25
$ws.Cells.Item(129, 14).Value = -13253.8531

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Cells.Item(5, 8).Value = 175
$ws.Cells.Item(5, 9).Value = 175
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 175
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -63

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Cells.Item(4, 8).Value = 175
$ws.Cells.Item(4, 9).Value = 175
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 175
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -60
# Row 129
$ws.Cells.Item(129, 8).Value = 49789.5
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 49789.5
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 49789.5
$ws.Cells.Item(129, 14).Value = -59789.5

$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Cells.Item(86, 8).Value = 5000
$ws.Cells.Item(86, 9).Value = 5000
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 5000
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = -3877
# Row 89
$ws.Cells.Item(89, 8).Value = 5000
$ws.Cells.Item(89, 9).Value = 5000
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 25000
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = -19384
# Row 121
$ws.Cells.Item(121, 8).Value = 22497.5
$ws.Cells.Item(121, 9).Value = 0
$ws.Cells.Item(121, 10).Value = 22497.5
$ws.Cells.Item(121, 11).Value = 0
$ws.Cells.Item(121, 12).Value = 22497.5
$ws.Cells.Item(121, 14).Value = -25117.5
# Row 135
$ws.Cells.Item(135, 8).Value = 87373
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 87373
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 87373
$ws.Cells.Item(135, 14).Value = -97513

$ws = $wb.Worksheets.Item("CUL")
# Row 62
$ws.Cells.Item(62, 8).Value = 4000
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 4000
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 12000
$ws.Cells.Item(62, 14).Value = -13372
# Row 63
$ws.Cells.Item(63, 8).Value = 4651.8335
$ws.Cells.Item(63, 9).Value = 3582.2
$ws.Cells.Item(63, 10).Value = 10000
$ws.Cells.Item(63, 11).Value = 10746.6
$ws.Cells.Item(63, 12).Value = 30000
$ws.Cells.Item(63, 13).Value = -9997.599999999999
$ws.Cells.Item(63, 14).Value = -31498
# Row 65
$ws.Cells.Item(65, 8).Value = 4000
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 4000
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 36000
$ws.Cells.Item(65, 14).Value = -42864
# Row 66
$ws.Cells.Item(66, 8).Value = 4651.8335
$ws.Cells.Item(66, 9).Value = 3582.2
$ws.Cells.Item(66, 10).Value = 10000
$ws.Cells.Item(66, 11).Value = 32239.8
$ws.Cells.Item(66, 12).Value = 90000
$ws.Cells.Item(66, 13).Value = -28495.8
$ws.Cells.Item(66, 14).Value = -97488
# Row 68
$ws.Cells.Item(68, 8).Value = 1408.7858
$ws.Cells.Item(68, 9).Value = 1036.6666
$ws.Cells.Item(68, 10).Value = 1687.875
$ws.Cells.Item(68, 11).Value = 3109.9998
$ws.Cells.Item(68, 12).Value = 5063.625
$ws.Cells.Item(68, 13).Value = -2298.9998
# Row 71
$ws.Cells.Item(71, 8).Value = 1408.7858
$ws.Cells.Item(71, 9).Value = 1036.6666
$ws.Cells.Item(71, 10).Value = 1687.875
$ws.Cells.Item(71, 11).Value = 9329.999400000001
$ws.Cells.Item(71, 12).Value = 15190.875
$ws.Cells.Item(71, 13).Value = -5273.999400000001
# Row 80
$ws.Cells.Item(80, 8).Value = 7099.875
$ws.Cells.Item(80, 9).Value = 10249.75
$ws.Cells.Item(80, 10).Value = 3950
$ws.Cells.Item(80, 11).Value = 30749.25
$ws.Cells.Item(80, 12).Value = 11850
$ws.Cells.Item(80, 13).Value = -29813.25
# Row 81
$ws.Cells.Item(81, 8).Value = 5005
$ws.Cells.Item(81, 9).Value = 2500
$ws.Cells.Item(81, 10).Value = 6257.5
$ws.Cells.Item(81, 11).Value = 7500
$ws.Cells.Item(81, 12).Value = 18772.5
$ws.Cells.Item(81, 13).Value = -6377
$ws.Cells.Item(81, 14).Value = -21018.5
# Row 82
$ws.Cells.Item(82, 8).Value = 12352
$ws.Cells.Item(82, 9).Value = 5000
$ws.Cells.Item(82, 10).Value = 13168.889
$ws.Cells.Item(82, 11).Value = 15000
$ws.Cells.Item(82, 12).Value = 39506.667
$ws.Cells.Item(82, 13).Value = -14594
$ws.Cells.Item(82, 14).Value = -40318.667
# Row 83
$ws.Cells.Item(83, 8).Value = 7099.875
$ws.Cells.Item(83, 9).Value = 10249.75
$ws.Cells.Item(83, 10).Value = 3950
$ws.Cells.Item(83, 11).Value = 92247.75
$ws.Cells.Item(83, 12).Value = 35550
$ws.Cells.Item(83, 13).Value = -87567.75
# Row 84
$ws.Cells.Item(84, 8).Value = 5005
$ws.Cells.Item(84, 9).Value = 2500
$ws.Cells.Item(84, 10).Value = 6257.5
$ws.Cells.Item(84, 11).Value = 22500
$ws.Cells.Item(84, 12).Value = 56317.5
$ws.Cells.Item(84, 13).Value = -16884
$ws.Cells.Item(84, 14).Value = -67549.5
# Row 85
$ws.Cells.Item(85, 8).Value = 12352
$ws.Cells.Item(85, 9).Value = 5000
$ws.Cells.Item(85, 10).Value = 13168.889
$ws.Cells.Item(85, 11).Value = 15000
$ws.Cells.Item(85, 12).Value = 39506.667
$ws.Cells.Item(85, 13).Value = -13596
$ws.Cells.Item(85, 14).Value = -42314.667
# Row 94
$ws.Cells.Item(94, 8).Value = 3810
$ws.Cells.Item(94, 9).Value = 3210
$ws.Cells.Item(94, 10).Value = 3960
$ws.Cells.Item(94, 11).Value = 9630
$ws.Cells.Item(94, 12).Value = 11880
$ws.Cells.Item(94, 13).Value = -8954
$ws.Cells.Item(94, 14).Value = -13232
# Row 97
$ws.Cells.Item(97, 8).Value = 576.6667
$ws.Cells.Item(97, 9).Value = 500
$ws.Cells.Item(97, 10).Value = 730
$ws.Cells.Item(97, 11).Value = 1500
$ws.Cells.Item(97, 12).Value = 2190
$ws.Cells.Item(97, 13).Value = -1004
$ws.Cells.Item(97, 14).Value = -3182
# Row 98
$ws.Cells.Item(98, 8).Value = 3825
$ws.Cells.Item(98, 9).Value = 3466.6667
$ws.Cells.Item(98, 10).Value = 4900
$ws.Cells.Item(98, 11).Value = 10400.0001
$ws.Cells.Item(98, 12).Value = 14700
$ws.Cells.Item(98, 13).Value = -8902.000100000001
$ws.Cells.Item(98, 14).Value = -17696
# Row 99
$ws.Cells.Item(99, 8).Value = 4246.2856
$ws.Cells.Item(99, 9).Value = 962.5
$ws.Cells.Item(99, 10).Value = 5559.8
$ws.Cells.Item(99, 11).Value = 2887.5
$ws.Cells.Item(99, 12).Value = 16679.4
$ws.Cells.Item(99, 13).Value = -641.5
$ws.Cells.Item(99, 14).Value = -21171.4
# Row 101
$ws.Cells.Item(101, 8).Value = 5451.84
$ws.Cells.Item(101, 9).Value = 0
$ws.Cells.Item(101, 10).Value = 5451.84
$ws.Cells.Item(101, 11).Value = 0
$ws.Cells.Item(101, 12).Value = 16355.52
$ws.Cells.Item(101, 14).Value = -21223.52
# Row 102
$ws.Cells.Item(102, 8).Value = 26750
$ws.Cells.Item(102, 9).Value = 50000
$ws.Cells.Item(102, 10).Value = 3500
$ws.Cells.Item(102, 11).Value = 150000
$ws.Cells.Item(102, 12).Value = 10500
$ws.Cells.Item(102, 13).Value = -147566
$ws.Cells.Item(102, 14).Value = -15368
# Row 103
$ws.Cells.Item(103, 8).Value = 2374.8
$ws.Cells.Item(103, 9).Value = 468.5
$ws.Cells.Item(103, 10).Value = 10000
$ws.Cells.Item(103, 11).Value = 1405.5
$ws.Cells.Item(103, 12).Value = 30000
$ws.Cells.Item(103, 13).Value = -526.5
# Row 121
$ws.Cells.Item(121, 8).Value = 72518.57000000001
$ws.Cells.Item(121, 9).Value = 1343.3334
$ws.Cells.Item(121, 10).Value = 125900
$ws.Cells.Item(121, 11).Value = 4030.0002
$ws.Cells.Item(121, 12).Value = 377700
$ws.Cells.Item(121, 13).Value = -2720.0002
$ws.Cells.Item(121, 14).Value = -380320
# Row 122
$ws.Cells.Item(122, 8).Value = 748.8461
$ws.Cells.Item(122, 9).Value = 498.3684
$ws.Cells.Item(122, 10).Value = 1428.7142
$ws.Cells.Item(122, 11).Value = 4485.3156
$ws.Cells.Item(122, 12).Value = 12858.4278
$ws.Cells.Item(122, 13).Value = -2035.3156
# Row 125
$ws.Cells.Item(125, 8).Value = 2991.25
$ws.Cells.Item(125, 9).Value = 1944
$ws.Cells.Item(125, 10).Value = 3467.2727
$ws.Cells.Item(125, 11).Value = 5832
$ws.Cells.Item(125, 12).Value = 10401.8181
$ws.Cells.Item(125, 13).Value = -912
$ws.Cells.Item(125, 14).Value = -20241.8181

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 3023.3333
$ws.Cells.Item(80, 9).Value = 3000
$ws.Cells.Item(80, 10).Value = 3052.5
$ws.Cells.Item(80, 11).Value = 3000
$ws.Cells.Item(80, 12).Value = 3052.5
$ws.Cells.Item(80, 13).Value = -2002
$ws.Cells.Item(80, 14).Value = -5048.5
# Row 83
$ws.Cells.Item(83, 8).Value = 3023.3333
$ws.Cells.Item(83, 9).Value = 3000
$ws.Cells.Item(83, 10).Value = 3052.5
$ws.Cells.Item(83, 11).Value = 15000
$ws.Cells.Item(83, 12).Value = 15262.5
$ws.Cells.Item(83, 13).Value = -10008
$ws.Cells.Item(83, 14).Value = -25246.5
# Row 132
$ws.Cells.Item(132, 8).Value = 2411.8
$ws.Cells.Item(132, 9).Value = 1620
$ws.Cells.Item(132, 10).Value = 3447.2307
$ws.Cells.Item(132, 11).Value = 4860
$ws.Cells.Item(132, 12).Value = 10341.6921
$ws.Cells.Item(132, 13).Value = -2330
$ws.Cells.Item(132, 14).Value = -15401.6921

$ws = $wb.Worksheets.Item("LTW")
# Row 129
$ws.Cells.Item(129, 8).Value = 60429
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 60429
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 60429
$ws.Cells.Item(129, 14).Value = -70429

$ws = $wb.Worksheets.Item("WVR")
# Row 123
$ws.Cells.Item(123, 8).Value = 54429
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 54429
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 54429
$ws.Cells.Item(123, 14).Value = -64229
# Row 129
$ws.Cells.Item(129, 8).Value = 24028.5
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 24028.5
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 24028.5
$ws.Cells.Item(129, 14).Value = -34028.5
